# "Added a None Option For Door" - the Door selection was set to None, so the
# DOOR line-item section (title, column headers, the "door exceeds glass area"
# row, and its trailing blank spacer row) no longer has anything to report and
# is removed from the report. The following FABRICATION section - and the
# totals below it - shift up to fill the gap. Because the report used to have
# an extra blank spacer row between the FABRICATION table and the SYSTEM TOTAL
# block, that now-redundant blank row collapses away too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/report tab.
$ws.Name = "Report"

# Remove the whole DOOR section: "DOOR" title (31), its Description/Part
# Number/Quantity/Price header row (32), the "Door size exceeds total glass
# area" data row (33), and the blank spacer row after it (34).
$ws.Rows("31:34").Delete()

# Collapse the now-duplicate blank spacer row that used to separate the
# FABRICATION table from the SYSTEM TOTAL block (old row 39, now row 35 after
# the shift above).
$ws.Rows("35:35").Delete()
